$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.626.59"
$ws.Range("E2").Value = "  +2.56%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.999.40"
$ws.Range("E3").Value = "  +6.08%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.33%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.32"
$ws.Range("E5").Value = "  +1.28%  "

# Row 6
$ws.Range("E6").Value = "  +0.17%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4670"
$ws.Range("E7").Value = "  +1.74%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3949"
$ws.Range("E8").Value = "  +1.49%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.76"
$ws.Range("E9").Value = "  +0.39%  "

# Row 10
$ws.Range("E10").Value = "  +1.54%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.005"
$ws.Range("E11").Value = "  +1.97%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.76"
$ws.Range("E12").Value = "  +4.42%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.040.55"
$ws.Range("E13").Value = "  +8.41%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.244"
$ws.Range("E14").Value = "  +3.26%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.876"
$ws.Range("E15").Value = "  +3.38%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07157"
$ws.Range("E16").Value = "  +3.19%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.02"
$ws.Range("E17").Value = "  +0.98%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.007"
$ws.Range("E18").Value = "  +0.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009998"
$ws.Range("E19").Value = "  +0.41%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.21"
$ws.Range("E20").Value = "  +1.25%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.730.75"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.540"
$ws.Range("E23").Value = "  +5.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.27"
$ws.Range("E24").Value = "  +2.82%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.266.78"
$ws.Range("E25").Value = "  +7.37%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.141"
$ws.Range("E26").Value = "  +2.68%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "158.17"
$ws.Range("E27").Value = "  +1.85%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.72"
$ws.Range("E28").Value = "  +2.06%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.999"
$ws.Range("E29").Value = "  +0.04%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.74"
$ws.Range("E30").Value = "  +2.72%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.967"
$ws.Range("E31").Value = "  +2.01%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09453"
$ws.Range("E32").Value = "  +1.28%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8950"
$ws.Range("E33").Value = "  -1.07%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.308"
$ws.Range("E34").Value = "  +0.56%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.348"
$ws.Range("E35").Value = "  +1.63%  "

# Row 36
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.185"
$ws.Range("E36").Value = "  -2.41%  "

# Row 37
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.000003767"
$ws.Range("E37").Value = "  +123.07%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05847"
$ws.Range("E38").Value = "  +1.50%  "

# Row 39
$ws.Range("E39").Value = "  -1.22%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02129"
$ws.Range("E40").Value = "  +2.78%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.938"
$ws.Range("E41").Value = "  +3.94%  "

# Row 42
$ws.Range("E42").Value = "  +0.11%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5775"
$ws.Range("E43").Value = "  +1.97%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1823"
$ws.Range("E44").Value = "  +3.21%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.900"
$ws.Range("E45").Value = "  +2.15%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.21"
$ws.Range("E46").Value = "  +2.46%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5405"
$ws.Range("E47").Value = "  +0.88%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.178"
$ws.Range("E48").Value = "  -3.11%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.652"
$ws.Range("E49").Value = "  +4.44%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07016"
$ws.Range("E50").Value = "  -0.32%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.876"
$ws.Range("E51").Value = "  +1.47%  "
